$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header J1: SIDEBAR_SUBMENU -> SUB_NAVBAR
$ws.Range("J1").Value = "SUB_NAVBAR"

# Update J2 value (was empty) -> Monitoring
$ws.Range("J2").Value = "Monitoring"

# Widen column J (10) to fit the new, longer content (closest attainable width to 18.5703125)
$ws.Columns.Item(10).ColumnWidth = 17.7

# Update selection / view to match the new active cell
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J2").Select()
